# Add new community summary rows (269-276) to the "Cells per mL" sheet,
# matching rows newly uploaded to the source spreadsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the existing date-formatted cell (D2) as the style donor so newly
# added date cells reuse the workbook's existing date style instead of
# minting a new one.
$dateDonor = $ws.Cells.Item(2, 4)

function Set-DataRow($Row, $Location, $Lat, $Lon, $DateSerial, $Counts, $Collector) {
    $ws.Cells.Item($Row, 1).Value = $Location
    $ws.Cells.Item($Row, 2).Value = $Lat
    $ws.Cells.Item($Row, 3).Value = $Lon

    $dCell = $ws.Cells.Item($Row, 4)
    $dCell.Value2 = $DateSerial
    $dateDonor.Copy()
    $dCell.PasteSpecial(-4122)

    foreach ($col in $Counts.Keys) {
        $ws.Cells.Item($Row, [int]$col).Value = $Counts[$col]
    }

    $ws.Cells.Item($Row, 16).Value = $Collector
}

# Row 269 - Wills Creek Price boat ramp
Set-DataRow 269 "Wills Creek Price boat ramp" `
    -34.289369999999998 138.01242999999999 45968 `
    @{6=56; 7=278; 8=611; 9=389; 10=667; 11=1778; 12=111; 13=833; 14=111; 15=4444} `
    "COLP5110"

# Row 270 - Garden Island
Set-DataRow 270 "Garden Island" `
    -34.804169999999999 138.53970000000001 45971 `
    @{8=3; 9=0; 11=42; 14=6; 15=51} `
    "HANC5232"

# Row 271 - Semaphore
Set-DataRow 271 "Semaphore" `
    -34.837600000000002 138.47788 45971 `
    @{8=2.5; 9=2.5; 11=180; 14=12.5; 15=195} `
    "HANC5232"

# Row 272 - Westlakes martin court
Set-DataRow 272 "Westlakes martin court" `
    -34.873939999999997 138.48820000000001 45971 `
    @{8=8; 9=0; 11=18.75; 14=16.25; 15=43.125} `
    "HANC5232"

# Row 273 - Hallet Cove (new location)
Set-DataRow 273 "Hallet Cove" `
    -35.078539999999997 138.49587 45971 `
    @{8=17; 9=10; 11=155; 14=0; 15=172} `
    "HANC5233"

# Row 274 - Kent Reserve Beach Victor Harbor (new location)
Set-DataRow 274 "Kent Reserve Beach Victor Harbor" `
    -35.563558 138.61283299999999 45970 `
    @{7=56; 8=333; 9=222; 10=56; 11=389; 13=167; 14=722; 15=1723} `
    "SEAS5211"

# Row 275 - Encounter Lake Tabernacle End
Set-DataRow 275 "Encounter Lake Tabernacle End" `
    -35.570839499999998 138.60133999999999 45972 `
    @{7=222; 9=0; 11=167; 13=111; 14=500; 15=1000} `
    "SEAS5211"

# Row 276 - Bluff Jetty
Set-DataRow 276 "Bluff Jetty" `
    -35.58858 138.60474600000001 45974 `
    @{6=111; 8=167; 9=111; 13=111; 14=444; 15=833} `
    "SEAS5211"

# Refresh the used-range selection to mirror where the author left off
# after the paste (D274:D276, the new date column values).
$ws.Range("D274:D276").Select()

Write-Host "Added rows 269-276"
